$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: model_V2
$ws.Range("A4").Value = "model_V2"
$ws.Range("B4").Value = 449
$ws.Range("C4").Value = 578
$ws.Range("D4").Value = 540
$ws.Range("E4").Value = 433
$ws.Range("F4").Formula = "=SUM(B4:E4)"
$ws.Range("G4").Formula = "=E4/F4"

# Copy styles from row 3 (E3, G3) to row 4 (E4, G4)
$ws.Range("E3").Copy() | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null
$ws.Range("G3").Copy() | Out-Null
$ws.Range("G4").PasteSpecial(-4122) | Out-Null

# New row 5: model_V3
$ws.Range("A5").Value = "model_V3"

# Match the final selection/active-cell state from the diff
$ws.Range("A6:F7").Select() | Out-Null
